# Generate Report for Handback
# Update timestamps on the "Overview", "zh-cn" and "de-de" sheets to
# reflect the latest handoff/handback generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-18 19:10:39"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-18 19:10:34"
$wsZhCn.Range("K2").Value = "2016-08-18 19:11:31"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-18 19:10:39"
$wsDeDe.Range("K2").Value = "2016-08-18 19:11:39"
